# Auto-generated edit script applying the Carbuncle_Profits (per-class Leve Profit) value updates.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 39132.125
$ws.Range("I21").Value = 37176.168
$ws.Range("K21").Value = 37176.168
$ws.Range("M21").Value = -36708.168
$ws.Range("H23").Value = 39132.125
$ws.Range("I23").Value = 37176.168
$ws.Range("K23").Value = 37176.168
$ws.Range("M23").Value = -36942.168
$ws.Range("H86").Value = 37488.184
$ws.Range("I86").Value = 1212.7142
$ws.Range("K86").Value = 1212.7142
$ws.Range("M86").Value = -89.71419999999989
$ws.Range("H89").Value = 37488.184
$ws.Range("I89").Value = 1212.7142
$ws.Range("K89").Value = 6063.571
$ws.Range("M89").Value = -447.5709999999999
$ws.Range("H106").Value = 2236.6875
$ws.Range("I106").Value = 1438.1818
$ws.Range("J106").Value = 3993.4
$ws.Range("K106").Value = 1438.1818
$ws.Range("L106").Value = 3993.4
$ws.Range("M106").Value = -807.1818000000001
$ws.Range("N106").Value = -5255.4
$ws.Range("H113").Value = 10959
$ws.Range("I113").Value = 2561.6667
$ws.Range("J113").Value = 13058.333
$ws.Range("K113").Value = 2561.6667
$ws.Range("L113").Value = 13058.333
$ws.Range("M113").Value = 692.3332999999998
$ws.Range("N113").Value = -19566.333
$ws.Range("H117").Value = 38249.5
$ws.Range("J117").Value = 38249.5
$ws.Range("L117").Value = 38249.5
$ws.Range("N117").Value = -47427.5
$ws.Range("H129").Value = 1116.75
$ws.Range("I129").Value = 481.625
$ws.Range("J129").Value = 1370.8
$ws.Range("K129").Value = 1444.875
$ws.Range("L129").Value = 4112.4
$ws.Range("M129").Value = 3555.125
$ws.Range("N129").Value = -14112.4
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()  # was -52140

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1812.5
$ws.Range("I2").Value = 1566.4
$ws.Range("J2").Value = 2058.6
$ws.Range("K2").Value = 1566.4
$ws.Range("L2").Value = 2058.6
$ws.Range("M2").Value = -1453.4
$ws.Range("N2").Value = -2284.6
$ws.Range("H74").Value = 2417.0557
$ws.Range("I74").Value = 2343.3125
$ws.Range("K74").Value = 2343.3125
$ws.Range("M74").Value = -1469.3125
$ws.Range("H77").Value = 2417.0557
$ws.Range("I77").Value = 2343.3125
$ws.Range("K77").Value = 11716.5625
$ws.Range("M77").Value = -7348.5625
$ws.Range("H116").Value = 1812.5
$ws.Range("I116").Value = 1566.4
$ws.Range("J116").Value = 2058.6
$ws.Range("K116").Value = 1566.4
$ws.Range("L116").Value = 2058.6
$ws.Range("M116").Value = 727.5999999999999
$ws.Range("N116").Value = -6646.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1812.5
$ws.Range("I3").Value = 1566.4
$ws.Range("J3").Value = 2058.6
$ws.Range("K3").Value = 1566.4
$ws.Range("L3").Value = 2058.6
$ws.Range("M3").Value = -1452.4
$ws.Range("N3").Value = -2286.6
$ws.Range("H86").Value = 2034.3334
$ws.Range("I86").Value = 2119.182
$ws.Range("J86").Value = 1901
$ws.Range("K86").Value = 2119.182
$ws.Range("L86").Value = 1901
$ws.Range("M86").Value = -996.1819999999998
$ws.Range("N86").Value = -4147
$ws.Range("H89").Value = 2034.3334
$ws.Range("I89").Value = 2119.182
$ws.Range("J89").Value = 1901
$ws.Range("K89").Value = 10595.91
$ws.Range("L89").Value = 9505
$ws.Range("M89").Value = -4979.91
$ws.Range("N89").Value = -20737
$ws.Range("H116").Value = 63155
$ws.Range("J116").Value = 63155
$ws.Range("L116").Value = 63155
$ws.Range("N116").Value = -72333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2912.0908
$ws.Range("I94").Value = 5750
$ws.Range("J94").Value = 2281.4443
$ws.Range("K94").Value = 5750
$ws.Range("L94").Value = 2281.4443
$ws.Range("M94").Value = -5299
$ws.Range("N94").Value = -3183.4443
$ws.Range("H132").Value = 2055.9688
$ws.Range("I132").Value = 1654.6
$ws.Range("J132").Value = 3489.4285
$ws.Range("K132").Value = 4963.799999999999
$ws.Range("L132").Value = 10468.2855
$ws.Range("M132").Value = -2433.799999999999
$ws.Range("N132").Value = -15528.2855

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 926.78125
$ws.Range("I68").Value = 714.4039
$ws.Range("J68").Value = 1177.7727
$ws.Range("K68").Value = 2143.2117
$ws.Range("L68").Value = 3533.3181
$ws.Range("M68").Value = -1332.2117
$ws.Range("N68").Value = -5155.3181
$ws.Range("H71").Value = 926.78125
$ws.Range("I71").Value = 714.4039
$ws.Range("J71").Value = 1177.7727
$ws.Range("K71").Value = 6429.6351
$ws.Range("L71").Value = 10599.9543
$ws.Range("M71").Value = -2373.6351
$ws.Range("N71").Value = -18711.9543
$ws.Range("H98").Value = 544.8333
$ws.Range("I98").Value = 589.6667
$ws.Range("J98").Value = 500
$ws.Range("K98").Value = 1769.0001
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = -271.0001
$ws.Range("N98").Value = -4496
$ws.Range("H107").Value = 709.7174
$ws.Range("I107").Value = 614.5
$ws.Range("K107").Value = 1843.5
$ws.Range("M107").Value = 76.5
$ws.Range("H122").Value = 372.65216
$ws.Range("I122").Value = 328.35715
$ws.Range("J122").Value = 441.55554
$ws.Range("K122").Value = 2955.21435
$ws.Range("L122").Value = 3973.99986
$ws.Range("M122").Value = -505.2143499999997
$ws.Range("N122").Value = -8873.99986
$ws.Range("H131").Value = 3548.9387
$ws.Range("J131").Value = 4026.1667
$ws.Range("L131").Value = 12078.5001
$ws.Range("N131").Value = -22158.5001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 53965
$ws.Range("J48").Value = 53965
$ws.Range("L48").Value = 53965
$ws.Range("N48").Value = -54935
$ws.Range("H113").Value = 1265.8334
$ws.Range("I113").Value = 1244.5454
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1244.5454
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 925.4546
$ws.Range("N113").Value = -5840

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 63715
$ws.Range("J36").Value = 63715
$ws.Range("L36").Value = 63715
$ws.Range("N36").Value = -64839
$ws.Range("H46").Value = 143914.28
$ws.Range("J46").Value = 1466.3334
$ws.Range("L46").Value = 1466.3334
$ws.Range("N46").Value = -1842.3334
$ws.Range("H61").Value = 2546.8572
$ws.Range("I61").Value = 2211.5
$ws.Range("J61").Value = 3620
$ws.Range("K61").Value = 2211.5
$ws.Range("L61").Value = 3620
$ws.Range("M61").Value = -2009.5
$ws.Range("N61").Value = -4024
$ws.Range("H113").Value = 2546.8572
$ws.Range("I113").Value = 2211.5
$ws.Range("J113").Value = 3620
$ws.Range("K113").Value = 2211.5
$ws.Range("L113").Value = 3620
$ws.Range("M113").Value = -41.5
$ws.Range("N113").Value = -7960
$ws.Range("H132").Value = 5296.514
$ws.Range("I132").Value = 5961.5884
$ws.Range("J132").Value = 4668.3887
$ws.Range("K132").Value = 17884.7652
$ws.Range("L132").Value = 14005.1661
$ws.Range("M132").Value = -15354.7652
$ws.Range("N132").Value = -19065.1661

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1368.75
$ws.Range("I113").Value = 1420.4546
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 4261.3638
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = -2091.3638
$ws.Range("N113").Value = -6740
$ws.Range("H132").Value = 2664.7
$ws.Range("I132").Value = 1975.6111
$ws.Range("K132").Value = 5926.8333
$ws.Range("M132").Value = -3396.8333
